$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Value = 111285339
$ws.Range("B29").Value = 94134
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 53
$ws.Range("F29").Value = "Vedtrappmossa"
$ws.Range("G29").Value = "Crossocalyx hellerianus"
$ws.Range("H29").Value = "(Nees ex Lindenb.) Meyl."
$ws.Range("I29").Value = ""
$ws.Range("J29").Value = ""
$ws.Range("Q29").Value = 463009.2170549285
$ws.Range("R29").Value = 6589537.119647364
$ws.Range("A30").Value = 111285077
$ws.Range("Q30").Value = 463029.9217482677
$ws.Range("R30").Value = 6589882.246201174
$ws.Range("A31").Value = 111285147
$ws.Range("B31").Value = 56398
$ws.Range("D31").Value = "NT"
$ws.Range("E31").Value = 100109
$ws.Range("F31").Value = "Tretåig hackspett"
$ws.Range("G31").Value = "Picoides tridactylus"
$ws.Range("H31").Value = "(Linnaeus, 1758)"
$ws.Range("J31").Value = ""
$ws.Range("M31").Value = "äldre spår"
$ws.Range("Q31").Value = 462998.224304053
$ws.Range("R31").Value = 6589870.859991242
$ws.Range("AF31").Value = ""
$ws.Range("A32").Value = 111285105
$ws.Range("B32").Value = 93067
$ws.Range("E32").Value = 2810
$ws.Range("F32").Value = "Västlig hakmossa"
$ws.Range("G32").Value = "Rhytidiadelphus loreus"
$ws.Range("H32").Value = "(Hedw.) Warnst."
$ws.Range("M32").Value = ""
$ws.Range("Q32").Value = 463000.5369220126
$ws.Range("R32").Value = 6589846.934736228
$ws.Range("Z32").Value = "00:00"
$ws.Range("AB32").Value = "00:00"
$ws.Range("A34").Value = 111285303
$ws.Range("B34").Value = 96348
$ws.Range("D34").Value = "VU"
$ws.Range("E34").Value = 220787
$ws.Range("F34").Value = "Knärot"
$ws.Range("G34").Value = "Goodyera repens"
$ws.Range("H34").Value = "(L.) R. Br."
$ws.Range("I34").Value = "'200"
$ws.Range("J34").Value = "stjälkar/strån/skott"
$ws.Range("Q34").Value = 463117.0211132796
$ws.Range("R34").Value = 6589724.244780275
$ws.Range("A35").Value = 111285012
$ws.Range("B35").Value = 4711
$ws.Range("D35").Value = "LC"
$ws.Range("E35").Value = 100299
$ws.Range("F35").Value = "Thomsons trägnagare"
$ws.Range("G35").Value = "Cacotemnus thomsoni"
$ws.Range("H35").Value = "(Kraatz, 1881)"
$ws.Range("J35").Value = ""
$ws.Range("M35").Value = ""
$ws.Range("Q35").Value = 463070.5957312917
$ws.Range("R35").Value = 6589821.837045968
$ws.Range("Z35").Value = "12:00"
$ws.Range("AB35").Value = "12:00"
$ws.Range("AF35").Value = ""
